# Update Betfair Back/Lay odds for 2025-12-26 with refreshed quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.94
$ws.Range("G2").Value = 1.97
$ws.Range("J2").Value = 4.3
$ws.Range("N2").Value = 6.6
$ws.Range("P2").Value = 2.9
$ws.Range("S2").Value = 2.22
$ws.Range("U2").Value = 2.82
$ws.Range("AH2").Value = 18.5

$ws.Range("F3").Value = 1.09
$ws.Range("J3").Value = 1.09
$ws.Range("V3").Value = 1.2

$ws.Range("G7").Value = 1.13
$ws.Range("P7").Value = 2.84
$ws.Range("Q7").Value = 1.23
$ws.Range("R7").Value = 2.14
$ws.Range("S7").Value = 1.55
$ws.Range("T7").Value = 2.04
$ws.Range("U7").Value = 1.44
$ws.Range("AN7").Value = 2.88

$ws.Range("F8").Value = 1.09
$ws.Range("G8").Value = 1.47
$ws.Range("K8").Value = 1000
$ws.Range("V8").Value = 1.02

$ws.Range("G9").Value = 3.9
$ws.Range("N9").Value = 2.08
$ws.Range("U9").Value = 2
$ws.Range("Y9").Value = 15.5
$ws.Range("AD9").Value = 17

$ws.Range("F10").Value = 2.32
$ws.Range("G10").Value = 2.44
$ws.Range("H10").Value = 2.92
$ws.Range("K10").Value = 4.2
$ws.Range("U10").Value = 1.04

$ws.Range("AJ12").Value = 980

$ws.Range("F14").Value = 2.1

$ws.Range("G16").Value = 2.22
$ws.Range("I16").Value = 4.3
$ws.Range("L16").Value = 1.29
$ws.Range("U16").Value = 2.1
$ws.Range("AK16").Value = 25

$ws.Range("O18").Value = 1.3
$ws.Range("T18").Value = 1.79
$ws.Range("AB18").Value = 11
$ws.Range("AK18").Value = 38

$ws.Range("J20").Value = 4.7
$ws.Range("P20").Value = 2.72
$ws.Range("R20").Value = 1.7
$ws.Range("S20").Value = 2.16
$ws.Range("V20").Value = 2.64

$ws.Range("Z21").Value = 11.5
$ws.Range("AG21").Value = 26
$ws.Range("AH21").Value = 24
$ws.Range("AI21").Value = 42

$ws.Range("F22").Value = 2.42

$ws.Range("J24").Value = 8
$ws.Range("W24").Value = 5.1
$ws.Range("X24").Value = 1000
$ws.Range("Z24").Value = 180
$ws.Range("AE24").Value = 190
$ws.Range("AF24").Value = 13
$ws.Range("AG24").Value = 13
$ws.Range("AH24").Value = 34
$ws.Range("AJ24").Value = 12
$ws.Range("AK24").Value = 13
$ws.Range("AO24").Value = 1000

$ws.Range("F25").Value = 2.62
$ws.Range("G25").Value = 2.64
$ws.Range("L25").Value = 1.29
$ws.Range("P25").Value = 2.38
$ws.Range("Q25").Value = 1.7
$ws.Range("R25").Value = 1.56
$ws.Range("V25").Value = 1.54
$ws.Range("W25").Value = 1.6
$ws.Range("X25").Value = 20
$ws.Range("Y25").Value = 14.5
$ws.Range("Z25").Value = 21
$ws.Range("AA25").Value = 42
$ws.Range("AB25").Value = 14.5
$ws.Range("AE25").Value = 26
$ws.Range("AF25").Value = 19.5
$ws.Range("AG25").Value = 11.5
$ws.Range("AJ25").Value = 38
$ws.Range("AK25").Value = 24
$ws.Range("AN25").Value = 16
$ws.Range("AO25").Value = 18.5
